$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Drop the four "Europe *" columns (old D:G); remaining columns shift left ---
$ws.Range("D1:G1").EntireColumn.Delete()

# --- Rewrite the two survey-question labels whose wording changed ---
# (built via formula + paste-values so the embedded line breaks do not trigger
#  an auto row-height resize; the final cells are plain static text, like the rest)
$ws.Range("A5").Formula = '="""Governments should actively cooperate to have all countries"&CHAR(10)&"converge in terms of GDP per capita by the end of the century"""'
$ws.Range("A6").Formula = '="Would support a global movement to tackle CC, tax millionaires,"&CHAR(10)&" and fund LICs (either petition, demonstrate, strike, or donate)*"'
$ws.Range("A5:A6").Copy()
$ws.Range("A5:A6").PasteSpecial(-4163)
$excel.CutCopyMode = $false

# --- Refresh the numeric data with the final computed values ---
$ws.Range("B2").Value = 0.693996120467564
$ws.Range("C2").Value = 0.410477795846558
$ws.Range("D2").Value = 0.670930964239829
$ws.Range("E2").Value = 0.744286257234234
$ws.Range("F2").Value = 0.669437974908923
$ws.Range("G2").Value = 0.817345461694808
$ws.Range("H2").Value = 0.848112605685568
$ws.Range("I2").Value = 0.725127730173656
$ws.Range("J2").Value = 0.770791010357976
$ws.Range("K2").Value = 0.422998865238901

$ws.Range("B3").Value = 0.641188431291775
$ws.Range("C3").Value = 0.369932629020835
$ws.Range("D3").Value = 0.632917148987785
$ws.Range("E3").Value = 0.670937221461187
$ws.Range("F3").Value = 0.457393328180411
$ws.Range("G3").Value = 0.816920760799125
$ws.Range("H3").Value = 0.837352904622364
$ws.Range("I3").Value = 0.676306813569887
$ws.Range("J3").Value = 0.697470165041275
$ws.Range("K3").Value = 0.414320739649134

$ws.Range("B4").Value = 0.680881448179833
$ws.Range("C4").Value = 0.616567982061628
$ws.Range("D4").Value = 0.743644347389163
$ws.Range("E4").Value = 0.814701212857562
$ws.Range("F4").Value = 0.757048871605567
$ws.Range("G4").Value = 0.713280127381035
$ws.Range("H4").Value = 0.703520370125625
$ws.Range("I4").Value = 0.671270631778761
$ws.Range("J4").Value = 0.776836935461012
$ws.Range("K4").Value = 0.425661149175785

$ws.Range("B5").Value = 0.704965329416964
$ws.Range("C5").Value = 0.501352982535347
$ws.Range("D5").Value = 0.7335602187152
$ws.Range("E5").Value = 0.725922165695082
$ws.Range("F5").Value = 0.685126203737904
$ws.Range("G5").Value = 0.930231790695484
$ws.Range("H5").Value = 0.938466265662205
$ws.Range("I5").Value = 0.688883535477258
$ws.Range("J5").Value = 0.68963486840272
$ws.Range("K5").Value = 0.436451458600216

$ws.Range("B6").Value = 0.675595447215337
$ws.Range("C6").Value = 0.523930159271177
$ws.Range("D6").Value = 0.433349195600366
$ws.Range("E6").Value = 0.696851480613757
$ws.Range("F6").Value = 0.583790255087382
$ws.Range("G6").Value = 0.727098526374066
$ws.Range("H6").Value = 0.741985444624183
$ws.Range("I6").Value = 0.641824096726743
$ws.Range("J6").Value = 0.834461320073758
$ws.Range("K6").Value = 0.474126518973143

$ws.Range("B7").Value = 0.682631646934764
$ws.Range("C7").Value = 0.485395103641793
$ws.Range("D7").Value = 0.421654543233796
$ws.Range("E7").Value = 0.721859577098009
$ws.Range("F7").Value = 0.567833043339243
$ws.Range("G7").ClearContents()
$ws.Range("H7").ClearContents()
$ws.Range("I7").Value = 0.594885214041605
$ws.Range("J7").Value = 0.89017740422893
$ws.Range("K7").Value = 0.463635219077665

$ws.Range("B8").Value = 0.451810364536854
$ws.Range("C8").Value = 0.300083084889478
$ws.Range("D8").ClearContents()
$ws.Range("E8").ClearContents()
$ws.Range("F8").ClearContents()
$ws.Range("G8").ClearContents()
$ws.Range("H8").ClearContents()
$ws.Range("I8").Value = 0.474757915976864
$ws.Range("J8").Value = 0.593250695193873
$ws.Range("K8").Value = 0.206703446618612

$ws.Range("B9").Value = 0.592122368373113
$ws.Range("C9").Value = 0.499097199497334
$ws.Range("D9").Value = 0.57447231505182
$ws.Range("E9").Value = 0.652952489891304
$ws.Range("F9").Value = 0.576305239322503
$ws.Range("G9").Value = 0.888820570273345
$ws.Range("H9").Value = 0.880416868630143
$ws.Range("I9").Value = 0.585354296646937
$ws.Range("J9").Value = 0.757501496513123
$ws.Range("K9").Value = 0.336242924079427
